$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "Datos actualizados..." timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 16 de Julio de 2020 a las 11:58"

# --- Swap country labels that moved places (Zambia <-> Eslovenia, Groenlandia <-> Islas Malvinas) ---
$ws.Range("A121").Value = "Eslovenia"
$ws.Range("A122").Value = "Zambia"

$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"

# --- Update numeric data for updated countries ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 3617474
$ws.Range("C4").Value = 647
$ws.Range("D4").Value = 1646675
$ws.Range("E4").Value = 1830639
$ws.Range("G4").Value = 16
$ws.Range("H4").Value = 140160

# Row 29: Indonesia
$ws.Range("B29").Value = 81668
$ws.Range("C29").Value = 1574
$ws.Range("D29").Value = 40345
$ws.Range("E29").Value = 37450
$ws.Range("G29").Value = 76
$ws.Range("H29").Value = 3873

# Row 35: Oman
$ws.Range("B35").Value = 62574
$ws.Range("C35").Value = 1327
$ws.Range("D35").Value = 40090
$ws.Range("E35").Value = 22194
$ws.Range("G35").Value = 9
$ws.Range("H35").Value = 290

# Row 36: Filipinas
$ws.Range("B36").Value = 61266
$ws.Range("C36").Value = 2498
$ws.Range("D36").Value = 21440
$ws.Range("E36").Value = 38183
$ws.Range("G36").Value = 29
$ws.Range("H36").Value = 1643

# Row 65: Marruecos
$ws.Range("B65").Value = 16424
$ws.Range("C65").Value = 162
$ws.Range("D65").Value = 13896
$ws.Range("E65").Value = 2268
$ws.Range("G65").Value = 1
$ws.Range("H65").Value = 260

# Row 80: Malasia
$ws.Range("B80").Value = 8737
$ws.Range("C80").Value = 3
$ws.Range("D80").Value = 8538
$ws.Range("E80").Value = 77

# Row 121: now Eslovenia (was Zambia's row, label just changed above)
$ws.Range("B121").Value = 1897
$ws.Range("C121").Value = 19
$ws.Range("D121").Value = 1522
$ws.Range("E121").Value = 264
$ws.Range("H121").Value = 111

# Row 122: now Zambia (was Eslovenia's row, label just changed above)
$ws.Range("B122").Value = 1895
$ws.Range("C122").Value = 0
$ws.Range("D122").Value = 1412
$ws.Range("E122").Value = 441
$ws.Range("H122").Value = 42

$wb.Save()
